$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24
$ws.Range("D24").Value = 'horseshoe crabs (Limulus polyphemus), red knots (Calidris canutus rufa)'
$ws.Range("E24").Value = 'Delaware Bay'
$ws.Range("F24").Value = 39
$ws.Range("G24").Value = -75
$ws.Range("H24").Value = 'Delaware Bay Horsewhoe crab fishery'
$ws.Range("I24").Value = 'Allowable Catch Adjustment, Monitoring Methodology, Species Interactions, Uncertainty'
$ws.Range("J24").Value = $true
$ws.Range("K24").Value = $true
$ws.Range("L24").Value = $true
$ws.Range("M24").Value = $true
$ws.Range("N24").Value = $true
$ws.Range("P24").Value = $true
$ws.Range("T24").Value = 'How to sustainably harvest [horseshoe crabs in Delaware Bay], including moratoria, given multispecies objectives?'
$ws.Range("U24").Value = 'Through consultations and public meetings with managers and stakeholders with additional input from scientific experts on sex ratio and carrying capacity.  While the objective function is reported, the full steps taken to produce it aren''t explicit.'
$ws.Range("V24").Value = 'MCDA'
$ws.Range("Y24").Value = 'Scientists'
$ws.Range("Z24").Value = 'Analysts, Decision Makers, Experts, Facilitators, Fishery, Government, Management, Scientists'
$ws.Range("AA24").Value = 'Analysts, Fishery, Management, Scientists'
$ws.Range("AC24").Value = 'Analysts, Fishery, Management'
$ws.Range("AE24").Value = 'Dynamic programming, Simulation modeling'
$ws.Range("AH24").Value = 'Study combines active adaptive management with management strategy evaluation.  Documentation of the process is not particularly detailed as presented here, but that stakeholder elicitiation occurred is documented. A value of information analysis was conducted in addition to the MSE.'

# Row 25
$ws.Range("D25").Value = 'NA, two generic simulated species (flathish and rockfish)'
$ws.Range("E25").Value = 'United States West Coast'
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = -125
$ws.Range("H25").Value = 'US west coast groundfish'
$ws.Range("I25").Value = 'Monitoring Methodology, Uncertainty'
$ws.Range("T25").Value = 'What are the impacts of monitoring method and stock assessment approachs associated with those methods to assessment performance?'
$ws.Range("U25").Value = 'Unspecified'
$ws.Range("W25").Value = 'Unknown'
$ws.Range("Y25").Value = 'Scientists'
$ws.Range("Z25").Value = 'Scientists'
$ws.Range("AB25").Value = 'Scientists'
$ws.Range("AD25").Value = 'Scientists'
$ws.Range("AE25").Value = 'Simulation modeling'
$ws.Range("AH25").Value = 'Simulation analysis seemingly unassociated with fishery problem development.  Monitoring methodology problem that did not take a VoI approach, or incorporate objectives beyond those related to estimation error.'

# Row 26
$ws.Range("D26").Value = 'eastern Baltic cod (Gadus morhua)'
$ws.Range("E26").Value = 'Baltic Sea'
$ws.Range("F26").Value = 59
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = 'Baltic Cod fishery'
$ws.Range("I26").Value = 'Environmental Conditions, Fishing Behavior'
$ws.Range("T26").Value = 'How will the recovery of Baltic cod be influenced by catch limit versus effort limit regulations, environmental conditions, and fleet behavior?'
$ws.Range("U26").Value = 'Unspecified.  Figures display more performance measures than described in the text.  Only those mentioned in the text are listed.'
$ws.Range("W26").Value = 'Unknown'
$ws.Range("Y26").Value = 'Unknown'
$ws.Range("Z26").Value = 'Scientists'
$ws.Range("AB26").Value = 'Scientists'
$ws.Range("AD26").Value = 'Management'
$ws.Range("AE26").Value = 'Simulation modeling'
$ws.Range("AH26").Value = 'This is a simulation analysis, but it may have been initiated due to management questioning the impact of recently implemented management changes.  The performance measures are not well documented, and the measures described do not correspond to the measures displayed in the results.'

# Row 27
$ws.Range("D27").Value = 'Southern Bluefin Tuna (Thunnus maccoyi)'
$ws.Range("E27").Value = 'Southern Ocean'
$ws.Range("F27").Value = -44
$ws.Range("G27").Value = 77
$ws.Range("H27").Value = 'Southern Bluefin Tuna Fishery'
$ws.Range("I27").Value = 'Uncertainty'
$ws.Range("J27").Value = $true
$ws.Range("L27").Value = $true
$ws.Range("M27").Value = $true
$ws.Range("N27").Value = $true
$ws.Range("O27").Value = $true
$ws.Range("T27").Value = 'What management proceedure will best meet the management objectives for southern bluefin tuna?'
$ws.Range("U27").Value = 'Through a series of workshops in which the management body provided objectives.  Only the SSB objective was explicitly established.  It seems that objective weights, or contraints for other objectives were attempted during the objective elicitation stage, which proved ineffective.'
$ws.Range("V27").Value = 'Visualization'
$ws.Range("Y27").Value = 'Management'
$ws.Range("Z27").Value = 'Decision Makers, Facilitators, Fishery, Government, Independent, Management, Scientists'
$ws.Range("AA27").Value = 'Management'
$ws.Range("AC27").Value = 'Scientists'
$ws.Range("AE27").Value = 'Simulation modeling'
$ws.Range("AH27").Value = 'This is a restrospective look back at an MSE process to derive lessons for fisheries management, and not exactly a documentation of a MSE itself.'

# Row 28
$ws.Range("D28").Value = 'sea lamprey (Petromyzon marinus)'
$ws.Range("E28").Value = 'Great Lakes, USA'
$ws.Range("F28").Value = 45
$ws.Range("G28").Value = -83
$ws.Range("H28").Value = 'integrated pest management (IPM) program implemented by the Great Lakes Fishery Commission (GLFC) to control invasive sea lamprey (Petromyzon marinus) in the Laurentian Great Lakes'
$ws.Range("I28").Value = 'Uncertainty'
$ws.Range("T28").Value = '"Managers implementing IPM must decide how to allocate resources among implementation of status quo control strategies, research, development, implementation of novel pest management strategies, and refinement of existing management strategies."'
$ws.Range("U28").Value = 'Unspecified'
$ws.Range("W28").Value = 'Unknown'
$ws.Range("Y28").Value = 'Unknown'
$ws.Range("Z28").Value = 'Management, Scientists'
$ws.Range("AB28").Value = 'Scientists'
$ws.Range("AD28").Value = 'Scientists'
$ws.Range("AE28").Value = 'Simulation modeling'
$ws.Range("AH28").Value = 'Simulaiton study without documentation of the decision making process.  Not the usual MSE focused on harvest management, but a case of aiming to cause mortality in a cost effective way.'

# Row 29
$ws.Range("D29").Value = 'Dungeness crab (Metacarcinus magister)'
$ws.Range("E29").Value = 'Hood Canal, Washington, USA'
$ws.Range("F29").Value = 48
$ws.Range("G29").Value = -123
$ws.Range("H29").Value = 'Hood Canal Dungeness crab fishery'
$ws.Range("I29").Value = 'Environmental Conditions, Implementation Uncertainty, Uncertainty'
$ws.Range("L29").Value = $true
$ws.Range("T29").Value = 'What Hood Canal Dungeness crab management strategies are robust to environmental and management uncertainty?'
$ws.Range("U29").Value = 'Co-managers’ objectives and interests. Several in-person and remote meetings with Washington State and tribal co-managers were organized to reveal key interests and considerations regarding harvest strategy for testing. The input from the co-managers gave a clearer conceptual understanding of the model content and assumptions.'
$ws.Range("W29").Value = 'Unknown'
$ws.Range("Y29").Value = 'Scientists'
$ws.Range("Z29").Value = 'Fishery, Government, Scientists'
$ws.Range("AA29").Value = 'Fishery, Government'
$ws.Range("AD29").Value = 'Unknown'
$ws.Range("AE29").Value = 'Expert elicitation, Simulation modeling'
$ws.Range("AH29").Value = 'Not a management strategy identification effort, but a consequence of environmental conditions study.'

# Row 30
$ws.Range("D30").Value = 'Southern rock lobster (Jasus edwardsii)'
$ws.Range("E30").Value = 'Great Australian Bight'
$ws.Range("F30").Value = -37
$ws.Range("G30").Value = 139
$ws.Range("H30").Value = 'Australian Southern Zone Southern rock lobster fishery'
$ws.Range("I30").Value = 'Uncertainty'
$ws.Range("K30").Value = $true
$ws.Range("M30").Value = $true
$ws.Range("O30").Value = $true
$ws.Range("S30").Value = $true
$ws.Range("T30").Value = 'Which of two proposed decision rules will best meet the objectives of the Australian Southern Zone Southern rock lobster fishery?'
$ws.Range("U30").Value = 'Unspecified'
$ws.Range("W30").Value = 'Mental Analysis'
$ws.Range("Y30").Value = 'Scientists'
$ws.Range("Z30").Value = 'Fishery, Government'
$ws.Range("AB30").Value = 'Fishery, Government'
$ws.Range("AC30").Value = 'Fishery, Government'
$ws.Range("AE30").Value = 'Simulation modeling'
$ws.Range("AH30").Value = 'Although the number of scenarios considered was not sufficient to conclude that all possible factors were considered, and the various scenarios were not weighted by the prior probabilities, the results of the analyses presented above assisted the selection of a DR for SZ rock lobster. The selected DR was the discrete rule. This selection was based on a number of factors, some based on results of the MSE and some based on other considerations. Specifically,
the performances of the two DRs were fairly similar with neither “dominating” the other. However, the discrete DR was seen to be more balanced (equally likely to increase as decrease the TACC) whereas it is more difficult to increase the TACC in the linear DR. Also, the discrete DR was proposed by the fishing industry and its acceptance increased stakeholder buy-in, a factor which cannot be over-emphasized during times when catch-rate [is in decline]'

# Row 31
$ws.Range("D31").Value = 'Yellow perch (Perca flavescens)'
$ws.Range("E31").Value = 'Lake Michigan, USA'
$ws.Range("F31").Value = 44
$ws.Range("G31").Value = -87
$ws.Range("H31").Value = 'Lake Michigan Yellow perch fishery'
$ws.Range("I31").Value = 'Uncertainty'
$ws.Range("J31").Value = $true
$ws.Range("K31").Value = $true
$ws.Range("L31").Value = $true
$ws.Range("M31").Value = $true
$ws.Range("Q31").Value = $true
$ws.Range("T31").Value = 'What harvest strategy should the Lake Michigan Committee and its constituent agencies adopt among different harvest policies to better meet their objectives for the fishery?'
$ws.Range("U31").Value = 'During the first workshop, participants developed a list of potential general management objectives and associated specific measures of performance during a facilitated brainstorming exercise.  Performance measures were refined through the series of workshops.'
$ws.Range("V31").Value = 'Unknown'
$ws.Range("Y31").Value = 'Government, Management, Scientists'
$ws.Range("Z31").Value = 'Decision Makers, Government, Management, Scientists'
$ws.Range("AA31").Value = 'Government, Management'
$ws.Range("AC31").Value = 'Government, Management'
$ws.Range("AE31").Value = 'Simulation modeling'

# Row 32
$ws.Range("AI32").Value = $false

# Row 33
$ws.Range("AI33").Value = $false

# Row 34
$ws.Range("AI34").Value = $false

# Row 35
$ws.Range("AI35").Value = $false

# Row 36
$ws.Range("AI36").Value = $false

# Row 37
$ws.Range("AI37").Value = $false

# Row 38
$ws.Range("AI38").Value = $false

# Row 39
$ws.Range("AI39").Value = $false

# Row 40
$ws.Range("AI40").Value = $false

# Row 41
$ws.Range("AI41").Value = $false
